$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 15 with the latest quotes (2025-09-19)
$ws.Cells.Item(15, 1).Value = 45919
$ws.Cells.Item(15, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item(15, 2).Value = "20,9437"
$ws.Cells.Item(15, 3).Value = "15,0727"
$ws.Cells.Item(15, 4).Value = "14,9476"
$ws.Cells.Item(15, 5).Value = "14,9476"
